$wb = $excel.ActiveWorkbook

# Sheet "展览" - update F3 (想去人数) from 994 to 995, and F4 (想去人数) from 492 to 493
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 995
$ws1.Range("F4").Value = 493

# Sheet "全部类型" - same updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 995
$ws4.Range("F4").Value = 493
